# species_types.xlsx edit script
# - Recapitalizes the two species-type labels ("evergreen" -> "Evergreen",
#   "deciduous" -> "Deciduous") while leaving "missing" as-is (this also
#   naturally reshuffles the shared-string table the same way the target
#   workbook's sharedStrings.xml/sheet1.xml does).
# - Re-selects cell U19 on the sheet (matches the saved cursor position).
# - Removes the pie chart's title and legend.
# - Turns on "show category name" for the data-label defaults as well as
#   for each individual pie slice's data label (category name + percentage,
#   no raw value) - matching the per-point <c:dLbl> blocks added in the diff.
# - Moves/resizes the chart to its new anchor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content -----------------------------------------------------
$ws.Range("A2").Value = "Evergreen"
$ws.Range("A3").Value = "Deciduous"
$ws.Range("A4").Value = "missing"

# --- Selection ----------------------------------------------------------
[void]$ws.Range("U19").Select()

# --- Chart ----------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

# Delete the chart title entirely (c:autoTitleDeleted val="1")
$chart.HasTitle = $false

# Delete the legend
$chart.HasLegend = $false

$series = $chart.SeriesCollection(1)

# Series-level data label defaults: show the category name too
$dlbls = $series.DataLabels()
$dlbls.ShowCategoryName = $true

# Per-point data labels (idx 0 = Evergreen slice, idx 1 = Deciduous slice):
# show category name + percentage, no raw value.
$pt1 = $series.Points().Item(1)
$pt2 = $series.Points().Item(2)

$dl1 = $pt1.DataLabel
$dl1.ShowValue = $false
$dl1.ShowCategoryName = $true
$dl1.ShowPercentage = $true

$dl2 = $pt2.DataLabel
$dl2.ShowValue = $false
$dl2.ShowCategoryName = $true
$dl2.ShowPercentage = $true

# --- Move / resize the chart -----------------------------------------
# New anchor: from col 10 (offset 1 EMU), row 7 (offset 95249 EMU)
#             to   col 17 (offset 457200 EMU), row 26 (offset 47624 EMU)
$co.Left = 582.9258599901575
$co.Top = 112.49992125984252
$co.Width = 445.0624212598425
$co.Height = 281.25
